$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.369.73'
$ws.Range("E2").Value = '  +0.35%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.008.56'
$ws.Range("E3").Value = '  -1.36%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '258.80'
$ws.Range("E5").Value = '  +4.30%  '

# Row 6
$ws.Range("E6").Value = '  -1.96%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.89'
$ws.Range("E8").Value = '  -6.32%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.384'
$ws.Range("E9").Value = '  -2.91%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0774'
$ws.Range("E10").Value = '  -5.06%  '

# Row 11
$ws.Range("E11").Value = '  -3.16%  '

# Row 12
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.306.87'
$ws.Range("E12").Value = '  -0.71%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.26'
$ws.Range("E13").Value = '  -6.73%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.20'
$ws.Range("E14").Value = '  -6.14%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.797'
$ws.Range("E15").Value = '  -7.88%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.22'
$ws.Range("E16").Value = '  -5.84%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.022.27'
$ws.Range("E17").Value = '  -0.01%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.263.81'
$ws.Range("E18").Value = '  +0.20%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.85'
$ws.Range("E19").Value = '  -1.42%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0838'
$ws.Range("E20").Value = '  -3.70%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '231.52'
$ws.Range("E21").Value = '  -0.05%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.11'
$ws.Range("E22").Value = '  -3.26%  '

# Row 23
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.04%  '

# Row 24
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.58'
$ws.Range("E24").Value = '  +2.37%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.34'
$ws.Range("E25").Value = '  -0.81%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.34'
$ws.Range("E26").Value = '  +0.24%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.94'
$ws.Range("E27").Value = '  -5.85%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.77'
$ws.Range("E28").Value = '  -0.88%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.128'
$ws.Range("E29").Value = '  -7.74%  '

# Row 30
$ws.Range("E30").Value = '  -4.03%  '

# Row 31
$ws.Range("E31").Value = '  -2.11%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0643'
$ws.Range("E32").Value = '  -4.50%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.57'
$ws.Range("E33").Value = '  -6.18%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.49'
$ws.Range("E34").Value = '  -1.36%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.36'
$ws.Range("E35").Value = '  -6.20%  '

# Row 36
$ws.Range("E36").Value = '  +0.48%  '

# Row 37
$ws.Range("E37").Value = '  +0.05%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.32'
$ws.Range("E38").Value = '  -4.16%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.43'
$ws.Range("E39").Value = '  -0.72%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.04'
$ws.Range("E40").Value = '  +3.19%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.19'
$ws.Range("E41").Value = '  -0.60%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0212'
$ws.Range("E42").Value = '  -1.90%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0927'
$ws.Range("E43").Value = '  -5.71%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.414.78'
$ws.Range("E44").Value = '  +1.58%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.75'
$ws.Range("E45").Value = '  -7.39%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.34'
$ws.Range("E46").Value = '  -3.63%  '

# Row 47
$ws.Range("E47").Value = '  -4.23%  '

# Row 48
$ws.Range("E48").Value = '  +1.66%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.02'
$ws.Range("E49").Value = '  -7.39%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.198.75'
$ws.Range("E50").Value = '  -0.86%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.94'
$ws.Range("E51").Value = '  -9.43%  '
